$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.55295333333333
$ws.Range("H2").Value = 37.65886
$ws.Range("I2").Value = 0.3363704472878066
$ws.Range("J2").Value = 0.3591006154861918
$ws.Range("M2").Value = 1.400501333333333
$ws.Range("N2").Value = 4.201504
$ws.Range("O2").Value = 0.00926314904242919
$ws.Range("P2").Value = 0.009687730200823723
$ws.Range("Q2").Value = 17.58042788060445
$ws.Range("R2").Value = 158.22385092544
$ws.Range("S2").Value = 0.003115849586695525
$ws.Range("T2").Value = 0.003478869877779968

$ws.Range("G3").Value = 12.55295333333333
$ws.Range("H3").Value = 37.65886
$ws.Range("I3").Value = 0.3363704472878066
$ws.Range("J3").Value = 0.3591006154861918
$ws.Range("O3").Value = 0.1405812059498714
$ws.Range("P3").Value = 0.1470248171880475
$ws.Range("Q3").Value = 266.8075123534867
$ws.Range("R3").Value = 2401.26761118138
$ws.Range("S3").Value = 0.0472873631256175
$ws.Range("T3").Value = 0.05279670234397268

$ws.Range("G4").Value = 12.55295333333333
$ws.Range("H4").Value = 37.65886
$ws.Range("I4").Value = 0.3363704472878066
$ws.Range("J4").Value = 0.3591006154861918
$ws.Range("M4").Value = 63.87756733333333
$ws.Range("N4").Value = 191.632702
$ws.Range("O4").Value = 0.4224968677952986
$ws.Range("P4").Value = 0.4418622271050682
$ws.Range("Q4").Value = 801.8521217821911
$ws.Range("R4").Value = 7216.669096039721
$ws.Range("S4").Value = 0.1421154603980019
$ws.Range("T4").Value = 0.1586729977135295

$ws.Range("G5").Value = 12.55295333333333
$ws.Range("H5").Value = 37.65886
$ws.Range("I5").Value = 0.3363704472878066
$ws.Range("J5").Value = 0.3591006154861918
$ws.Range("M5").Value = 19.878555
$ws.Range("N5").Value = 39.75711
$ws.Range("O5").Value = 0.1314800731212866
$ws.Range("P5").Value = 0.0916710195312133
$ws.Range("Q5").Value = 249.5345732491
$ws.Range("R5").Value = 1497.2074394946
$ws.Range("S5").Value = 0.04422601100524069
$ws.Range("T5").Value = 0.03291911953590541

$ws.Range("G6").Value = 12.55295333333333
$ws.Range("H6").Value = 37.65886
$ws.Range("I6").Value = 0.3363704472878066
$ws.Range("J6").Value = 0.3591006154861918
$ws.Range("M6").Value = 44.77944466666667
$ws.Range("N6").Value = 134.338334
$ws.Range("O6").Value = 0.2961787040911142
$ws.Range("P6").Value = 0.3097542059748472
$ws.Range("Q6").Value = 562.114279193249
$ws.Range("R6").Value = 5059.028512739241
$ws.Range("S6").Value = 0.09962576317225103
$ws.Range("T6").Value = 0.1112329260150042

$ws.Range("H7").Value = 39.447015
$ws.Range("I7").Value = 0.3523423194360853
$ws.Range("J7").Value = 0.3761517838190811
$ws.Range("M7").Value = 1.400501333333333
$ws.Range("N7").Value = 4.201504
$ws.Range("O7").Value = 0.00926314904242919
$ws.Range("P7").Value = 0.009687730200823723
$ws.Range("Q7").Value = 18.41519903450667
$ws.Range("R7").Value = 165.73679131056
$ws.Range("S7").Value = 0.003263799418891654
$ws.Range("T7").Value = 0.003644056996197828

$ws.Range("H8").Value = 39.447015
$ws.Range("I8").Value = 0.3523423194360853
$ws.Range("J8").Value = 0.3761517838190811
$ws.Range("O8").Value = 0.1405812059498714
$ws.Range("P8").Value = 0.1470248171880475
$ws.Range("Q8").Value = 279.476328861805
$ws.Range("R8").Value = 2515.286959756245
$ws.Range("S8").Value = 0.04953270817349969
$ws.Range("T8").Value = 0.05530364725095835

$ws.Range("H9").Value = 39.447015
$ws.Range("I9").Value = 0.3523423194360853
$ws.Range("J9").Value = 0.3761517838190811
$ws.Range("M9").Value = 63.87756733333333
$ws.Range("N9").Value = 191.632702
$ws.Range("O9").Value = 0.4224968677952986
$ws.Range("P9").Value = 0.4418622271050682
$ws.Range("Q9").Value = 839.9264522538367
$ws.Range("R9").Value = 7559.33807028453
$ws.Range("S9").Value = 0.1488635263534766
$ws.Range("T9").Value = 0.1662072649278434

$ws.Range("H10").Value = 39.447015
$ws.Range("I10").Value = 0.3523423194360853
$ws.Range("J10").Value = 0.3761517838190811
$ws.Range("M10").Value = 19.878555
$ws.Range("N10").Value = 39.75711
$ws.Range("O10").Value = 0.1314800731212866
$ws.Range("P10").Value = 0.0916710195312133
$ws.Range("Q10").Value = 261.383219087775
$ws.Range("R10").Value = 1568.29931452665
$ws.Range("S10").Value = 0.04632599392318021
$ws.Range("T10").Value = 0.03448221752117971

$ws.Range("H11").Value = 39.447015
$ws.Range("I11").Value = 0.3523423194360853
$ws.Range("J11").Value = 0.3761517838190811
$ws.Range("M11").Value = 44.77944466666667
$ws.Range("N11").Value = 134.338334
$ws.Range("O11").Value = 0.2961787040911142
$ws.Range("P11").Value = 0.3097542059748472
$ws.Range("Q11").Value = 588.8051418192234
$ws.Range("R11").Value = 5299.24627637301
$ws.Range("S11").Value = 0.1043562915670372
$ws.Range("T11").Value = 0.1165145971229018

$ws.Range("G12").Value = 2.132104
$ws.Range("H12").Value = 6.396312
$ws.Range("I12").Value = 0.05713211521624299
$ws.Range("J12").Value = 0.06099280689967021
$ws.Range("M12").Value = 1.400501333333333
$ws.Range("N12").Value = 4.201504
$ws.Range("O12").Value = 0.00926314904242919
$ws.Range("P12").Value = 0.009687730200823723
$ws.Range("Q12").Value = 2.986014494805334
$ws.Range("R12").Value = 26.874130453248
$ws.Range("S12").Value = 0.0005292232983572955
$ws.Range("T12").Value = 0.0005908818574349446

$ws.Range("G13").Value = 2.132104
$ws.Range("H13").Value = 6.396312
$ws.Range("I13").Value = 0.05713211521624299
$ws.Range("J13").Value = 0.06099280689967021
$ws.Range("O13").Value = 0.1405812059498714
$ws.Range("P13").Value = 0.1470248171880475
$ws.Range("Q13").Value = 45.316934526344
$ws.Range("R13").Value = 407.852410737096
$ws.Range("S13").Value = 0.008031701655566437
$ws.Range("T13").Value = 0.008967456284209892

$ws.Range("G14").Value = 2.132104
$ws.Range("H14").Value = 6.396312
$ws.Range("I14").Value = 0.05713211521624299
$ws.Range("J14").Value = 0.06099280689967021
$ws.Range("M14").Value = 63.87756733333333
$ws.Range("N14").Value = 191.632702
$ws.Range("O14").Value = 0.4224968677952986
$ws.Range("P14").Value = 0.4418622271050682
$ws.Range("Q14").Value = 136.1936168216693
$ws.Range("R14").Value = 1225.742551395024
$ws.Range("S14").Value = 0.02413813972938278
$ws.Range("T14").Value = 0.02695041749407765

$ws.Range("G15").Value = 2.132104
$ws.Range("H15").Value = 6.396312
$ws.Range("I15").Value = 0.05713211521624299
$ws.Range("J15").Value = 0.06099280689967021
$ws.Range("M15").Value = 19.878555
$ws.Range("N15").Value = 39.75711
$ws.Range("O15").Value = 0.1314800731212866
$ws.Range("P15").Value = 0.0916710195312133
$ws.Range("Q15").Value = 42.38314662972
$ws.Range("R15").Value = 254.29887977832
$ws.Range("S15").Value = 0.007511734686205398
$ws.Range("T15").Value = 0.005591272792563189

$ws.Range("G16").Value = 2.132104
$ws.Range("H16").Value = 6.396312
$ws.Range("I16").Value = 0.05713211521624299
$ws.Range("J16").Value = 0.06099280689967021
$ws.Range("M16").Value = 44.77944466666667
$ws.Range("N16").Value = 134.338334
$ws.Range("O16").Value = 0.2961787040911142
$ws.Range("P16").Value = 0.3097542059748472
$ws.Range("Q16").Value = 95.47443309157867
$ws.Range("R16").Value = 859.269897824208
$ws.Range("S16").Value = 0.01692131584673108
$ws.Range("T16").Value = 0.01889277847138452

$ws.Range("G17").Value = 7.086566
$ws.Range("H17").Value = 14.173132
$ws.Range("I17").Value = 0.1898924748509033
$ws.Range("J17").Value = 0.135149614846733
$ws.Range("M17").Value = 1.400501333333333
$ws.Range("N17").Value = 4.201504
$ws.Range("O17").Value = 0.00926314904242919
$ws.Range("P17").Value = 0.009687730200823723
$ws.Range("Q17").Value = 9.924745131754667
$ws.Range("R17").Value = 59.548470790528
$ws.Range("S17").Value = 0.001759002296579654
$ws.Range("T17").Value = 0.00130929300538039

$ws.Range("G18").Value = 7.086566
$ws.Range("H18").Value = 14.173132
$ws.Range("I18").Value = 0.1898924748509033
$ws.Range("J18").Value = 0.135149614846733
$ws.Range("O18").Value = 0.1405812059498714
$ws.Range("P18").Value = 0.1470248171880475
$ws.Range("Q18").Value = 150.621849327526
$ws.Range("R18").Value = 903.7310959651561
$ws.Range("S18").Value = 0.02669531311534561
$ws.Range("T18").Value = 0.01987034741587595

$ws.Range("G19").Value = 7.086566
$ws.Range("H19").Value = 14.173132
$ws.Range("I19").Value = 0.1898924748509033
$ws.Range("J19").Value = 0.135149614846733
$ws.Range("M19").Value = 63.87756733333333
$ws.Range("N19").Value = 191.632702
$ws.Range("O19").Value = 0.4224968677952986
$ws.Range("P19").Value = 0.4418622271050682
$ws.Range("Q19").Value = 452.6725968271107
$ws.Range("R19").Value = 2716.035580962664
$ws.Range("S19").Value = 0.08022897584240414
$ws.Range("T19").Value = 0.05971750980856966

$ws.Range("G20").Value = 7.086566
$ws.Range("H20").Value = 14.173132
$ws.Range("I20").Value = 0.1898924748509033
$ws.Range("J20").Value = 0.135149614846733
$ws.Range("M20").Value = 19.878555
$ws.Range("N20").Value = 39.75711
$ws.Range("O20").Value = 0.1314800731212866
$ws.Range("P20").Value = 0.0916710195312133
$ws.Range("Q20").Value = 140.87069199213
$ws.Range("R20").Value = 563.48276796852
$ws.Range("S20").Value = 0.02496707647857883
$ws.Range("T20").Value = 0.01238930298225082

$ws.Range("G21").Value = 7.086566
$ws.Range("H21").Value = 14.173132
$ws.Range("I21").Value = 0.1898924748509033
$ws.Range("J21").Value = 0.135149614846733
$ws.Range("M21").Value = 44.77944466666667
$ws.Range("N21").Value = 134.338334
$ws.Range("O21").Value = 0.2961787040911142
$ws.Range("P21").Value = 0.3097542059748472
$ws.Range("Q21").Value = 317.3324900736814
$ws.Range("R21").Value = 1903.994940442088
$ws.Range("S21").Value = 0.05624210711799503
$ws.Range("T21").Value = 0.04186316163465621

$ws.Range("G22").Value = 2.398207
$ws.Range("H22").Value = 7.194621
$ws.Range("I22").Value = 0.06426264320896187
$ws.Range("J22").Value = 0.06860517894832399
$ws.Range("M22").Value = 1.400501333333333
$ws.Range("N22").Value = 4.201504
$ws.Range("O22").Value = 0.00926314904242919
$ws.Range("P22").Value = 0.009687730200823723
$ws.Range("Q22").Value = 3.358692101109333
$ws.Range("R22").Value = 30.228228909984
$ws.Range("S22").Value = 0.0005952744419050639
$ws.Range("T22").Value = 0.0006646284640305942

$ws.Range("G23").Value = 2.398207
$ws.Range("H23").Value = 7.194621
$ws.Range("I23").Value = 0.06426264320896187
$ws.Range("J23").Value = 0.06860517894832399
$ws.Range("O23").Value = 0.1405812059498714
$ws.Range("P23").Value = 0.1470248171880475
$ws.Range("Q23").Value = 50.97283697212699
$ws.Range("R23").Value = 458.755532749143
$ws.Range("S23").Value = 0.009034119879842173
$ws.Range("T23").Value = 0.01008666389303062

$ws.Range("G24").Value = 2.398207
$ws.Range("H24").Value = 7.194621
$ws.Range("I24").Value = 0.06426264320896187
$ws.Range("J24").Value = 0.06860517894832399
$ws.Range("M24").Value = 63.87756733333333
$ws.Range("N24").Value = 191.632702
$ws.Range("O24").Value = 0.4224968677952986
$ws.Range("P24").Value = 0.4418622271050682
$ws.Range("Q24").Value = 153.1916291217713
$ws.Range("R24").Value = 1378.724662095942
$ws.Range("S24").Value = 0.02715076547203321
$ws.Range("T24").Value = 0.03031403716104818

$ws.Range("G25").Value = 2.398207
$ws.Range("H25").Value = 7.194621
$ws.Range("I25").Value = 0.06426264320896187
$ws.Range("J25").Value = 0.06860517894832399
$ws.Range("M25").Value = 19.878555
$ws.Range("N25").Value = 39.75711
$ws.Range("O25").Value = 0.1314800731212866
$ws.Range("P25").Value = 0.0916710195312133
$ws.Range("Q25").Value = 47.67288975088499
$ws.Range("R25").Value = 286.0373385053099
$ws.Range("S25").Value = 0.008449257028081456
$ws.Range("T25").Value = 0.006289106699314193

$ws.Range("G26").Value = 2.398207
$ws.Range("H26").Value = 7.194621
$ws.Range("I26").Value = 0.06426264320896187
$ws.Range("J26").Value = 0.06860517894832399
$ws.Range("M26").Value = 44.77944466666667
$ws.Range("N26").Value = 134.338334
$ws.Range("O26").Value = 0.2961787040911142
$ws.Range("P26").Value = 0.3097542059748472
$ws.Range("Q26").Value = 107.3903776557127
$ws.Range("R26").Value = 966.513398901414
$ws.Range("S26").Value = 0.01903322638709997
$ws.Range("T26").Value = 0.0212507427309004
